$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04155
$ws.Range("H2").Value = 0.12465
$ws.Range("I2").Value = 0.0001466168179836329
$ws.Range("J2").Value = 0.0001466168179836329
$ws.Range("M2").Value = 0.029424
$ws.Range("N2").Value = 0.08827199999999999
$ws.Range("O2").Value = 0.1473063425232919
$ws.Range("P2").Value = 0.1473063425232919
$ws.Range("Q2").Value = 0.0012225672
$ws.Range("R2").Value = 0.0110031048
$ws.Range("S2").Value = 0.00002159758720957217
$ws.Range("T2").Value = 0.00002159758720957217

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04155
$ws.Range("H3").Value = 0.12465
$ws.Range("I3").Value = 0.0001466168179836329
$ws.Range("J3").Value = 0.0001466168179836329
$ws.Range("O3").Value = 0.852693657476708
$ws.Range("P3").Value = 0.852693657476708
$ws.Range("Q3").Value = 0.00707692065
$ws.Range("R3").Value = 0.06369228585
$ws.Range("S3").Value = 0.0001250192307740607
$ws.Range("T3").Value = 0.0001250192307740607

# Row 4
$ws.Range("I4").Value = 0.9992428949822291
$ws.Range("J4").Value = 0.9992428949822291
$ws.Range("M4").Value = 0.029424
$ws.Range("N4").Value = 0.08827199999999999
$ws.Range("O4").Value = 0.1473063425232919
$ws.Range("P4").Value = 0.1473063425232919
$ws.Range("Q4").Value = 8.332206393775998
$ws.Range("R4").Value = 74.98985754398399
$ws.Range("S4").Value = 0.1471948161522181
$ws.Range("T4").Value = 0.1471948161522181

# Row 5
$ws.Range("I5").Value = 0.9992428949822291
$ws.Range("J5").Value = 0.9992428949822291
$ws.Range("O5").Value = 0.852693657476708
$ws.Range("P5").Value = 0.852693657476708
$ws.Range("S5").Value = 0.852048078830011
$ws.Range("T5").Value = 0.852048078830011

# Row 6
$ws.Range("I6").Value = 0.0006104881997874136
$ws.Range("J6").Value = 0.0006104881997874135
$ws.Range("M6").Value = 0.029424
$ws.Range("N6").Value = 0.08827199999999999
$ws.Range("O6").Value = 0.1473063425232919
$ws.Range("P6").Value = 0.1473063425232919
$ws.Range("Q6").Value = 0.005090567775999998
$ws.Range("R6").Value = 0.04581510998399999
$ws.Range("S6").Value = 0.00008992878386431262
$ws.Range("T6").Value = 0.00008992878386431261

# Row 7
$ws.Range("I7").Value = 0.0006104881997874136
$ws.Range("J7").Value = 0.0006104881997874135
$ws.Range("O7").Value = 0.852693657476708
$ws.Range("P7").Value = 0.852693657476708
$ws.Range("S7").Value = 0.000520559415923101
$ws.Range("T7").Value = 0.0005205594159231008
